$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# --- Style conversions: cells that change from "N/A" text to numeric values ---
# (style 15 = #,##0 integer format; style 16 = #,##0.0;"-"#,##0.0 percent-like format)
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C29").NumberFormat = "#,##0"

# --- Data value updates, rows 14-29 ---
# Row 14
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("M14").Value = -20
$ws.Range("N14").Value = -81.25

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 28
$ws.Range("J15").Value = 23
$ws.Range("K15").Value = 21.739130434782
$ws.Range("L15").Value = 27.272727272727
$ws.Range("N15").Value = -47.169811320754

# Row 16
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 53
$ws.Range("G16").Value = 49
$ws.Range("H16").Value = 8.163265306122
$ws.Range("I16").Value = 322
$ws.Range("J16").Value = 389
$ws.Range("K16").Value = -17.223650385604
$ws.Range("L16").Value = 30.894308943089
$ws.Range("M16").Value = -6.122448979591
$ws.Range("N16").Value = -76.632801161103

# Row 17
$ws.Range("C17").Value = 25
$ws.Range("D17").Value = 27
$ws.Range("E17").Value = -7.407407407407
$ws.Range("F17").Value = 92
$ws.Range("G17").Value = 87
$ws.Range("H17").Value = 5.747126436781
$ws.Range("I17").Value = 512
$ws.Range("J17").Value = 523
$ws.Range("K17").Value = -2.103250478011
$ws.Range("L17").Value = 31.282051282051
$ws.Range("M17").Value = 45.042492917847
$ws.Range("N17").Value = -27.785613540197

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -23.076923076923
$ws.Range("I18").Value = 185
$ws.Range("J18").Value = 199
$ws.Range("K18").Value = -7.035175879396
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -12.735849056603
$ws.Range("N18").Value = -78.881278538812

# Row 19
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 106
$ws.Range("H19").Value = -4.716981132075
$ws.Range("I19").Value = 531
$ws.Range("J19").Value = 635
$ws.Range("K19").Value = -16.377952755905
$ws.Range("L19").Value = 12.738853503184
$ws.Range("M19").Value = 68.037974683544
$ws.Range("N19").Value = 27.951807228915

# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -36.363636363636
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 3.571428571428
$ws.Range("I20").Value = 256
$ws.Range("J20").Value = 251
$ws.Range("K20").Value = 1.992031872509
$ws.Range("L20").Value = 8.936170212765
$ws.Range("M20").Value = 77.777777777777
$ws.Range("N20").Value = -80.664652567975

# Row 21
$ws.Range("C21").Value = 76
$ws.Range("D21").Value = 80
$ws.Range("E21").Value = -5
$ws.Range("F21").Value = 300
$ws.Range("H21").Value = 0.334448160535
$ws.Range("I21").Value = 1846
$ws.Range("J21").Value = 2028
$ws.Range("K21").Value = -8.974358974358
$ws.Range("L21").Value = 18.485237483953
$ws.Range("M21").Value = 30.829199149539
$ws.Range("N21").Value = -61.693297364598

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = -54.166666666666
$ws.Range("L22").Value = 46.666666666666

# Row 23
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 40
$ws.Range("F23").Value = 37
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = 68.181818181818
$ws.Range("I23").Value = 180
$ws.Range("J23").Value = 167
$ws.Range("K23").Value = 7.784431137724
$ws.Range("L23").Value = 29.496402877697
$ws.Range("M23").Value = 102.247191011236

# Row 24
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = 14.035087719298
$ws.Range("F24").Value = 205
$ws.Range("G24").Value = 214
$ws.Range("H24").Value = -4.205607476635
$ws.Range("I24").Value = 1138
$ws.Range("J24").Value = 1399
$ws.Range("K24").Value = -18.656182987848
$ws.Range("L24").Value = 6.654170571696
$ws.Range("M24").Value = 60.056258790436

# Row 25
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 119
$ws.Range("G25").Value = 137
$ws.Range("H25").Value = -13.138686131386
$ws.Range("I25").Value = 598
$ws.Range("J25").Value = 658
$ws.Range("K25").Value = -9.118541033434
$ws.Range("L25").Value = 36.529680365296
$ws.Range("M25").Value = -31.343283582089

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 9
$ws.Range("H26").Value = 125
$ws.Range("I26").Value = 43
$ws.Range("J26").Value = 36
$ws.Range("K26").Value = 19.444444444444
$ws.Range("L26").Value = -4.444444444444

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = -40
$ws.Range("F27").Value = 12
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 58
$ws.Range("J27").Value = 58
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5.454545454545

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -70
$ws.Range("I28").Value = 38
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = -11.627906976744
$ws.Range("L28").Value = 11.764705882352
$ws.Range("M28").Value = -9.523809523809
$ws.Range("N28").Value = -80.104712041884

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -50
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = -57.142857142857
$ws.Range("I29").Value = 32
$ws.Range("J29").Value = 37
$ws.Range("K29").Value = -13.513513513513
$ws.Range("L29").Value = -3.030303030303
$ws.Range("M29").Value = -8.571428571428
$ws.Range("N29").Value = -81.609195402298

